$wb = $excel.ActiveWorkbook

# --- Sheet references (by position: 1=Notes, 2=Лист1, 3=Лист3) ---
$wsNotes = $wb.Worksheets.Item(1)
$wsChrom = $wb.Worksheets.Item(2)

# --- 1) Rename "Лист1" -> "Chrom to dia" ---
$wsChrom.Name = "Chrom to dia"

# --- 2) Convert the column-B "fill" formulas on Notes into shared formulas
#        (matches the pattern: B2:B10 fills from B9+12, B18:B32 fills from B11-12) ---
$wsNotes.Range("B2:B10").Formula = "=B9+12"
$wsNotes.Range("B18:B32").Formula = "=B11-12"

# --- 3) Populate "Chrom to dia" sheet with the new analysis data ---
$wsChrom.Range("A1").Value = "Chrom"
$startVal = 86
for ($r = 2; $r -le 42; $r++) {
    $wsChrom.Cells.Item($r, 1).Value = $startVal - ($r - 2)
}

# --- 4) Make "Chrom to dia" the active sheet/tab, with C3 selected ---
$wsChrom.Activate()
$wsChrom.Range("C3").Select()
